$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert 3 blank rows before the old "blank separator" rows (old rows 9-10)
#    so the table grows from 12 to 15 rows, pushing the totals/pages rows
#    down and widening the SUM() ranges automatically.
# ---------------------------------------------------------------------------
$ws.Rows("9:11").Insert()
$ws.Rows("9:11").RowHeight = 15

# ---------------------------------------------------------------------------
# 2. Move Phillis/David comment from the cell that used to hold
#    "Happenin's 'Round the House" (old A5) to its new location (A8).
# ---------------------------------------------------------------------------
$commentText = $ws.Range("A5").Comment.Text()
$ws.Range("A5").Comment.Delete()
$ws.Range("A8").AddComment($commentText) | Out-Null

# ---------------------------------------------------------------------------
# 3. Rewrite the whole article table (rows 2-11) with the updated data.
# ---------------------------------------------------------------------------
$data = @(
    @("Silents",                     "Bill Blowers",   500, 1, 448, 9),
    @("Critters",                    "Rachael K.",     500, 3, 520, 4),
    @("Goin Country",                "Roger Basham",   $null, 2, 14, 2),
    @("Hart Books",                  "Jennie????",     500, 2, 479, 3),
    @("Hart's Other Home",           "Bill West",      200, 3, 106, 3),
    @("Trading Post",                "Roger Basham",   500, 1, 222, 1),
    @("Happenin's 'Round the House", "Margi Bertram",  500, 1, 146, 1),
    @("Meet a board member",         "Tim Murphy",     500, 2, 550, 2),
    @("Pow Wow & other Flyers",      "Bob Hoke",       600, 2, 600, 2),
    @("Back Page",                   "me",             600, 0, 600, $null)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $r++
}

# Row 11's G cell must stay empty (no value at all, like the old template rows).
$ws.Cells.Item(11, 7).Value = $null

# ---------------------------------------------------------------------------
# 4. The hidden helper shape behind the legacy comment box is two-cell
#    anchored all the way down around row 80; growing the table by three
#    rows pushes its bottom anchor down to row 83 as well.
# ---------------------------------------------------------------------------
$shape = $ws.Shapes.Item(1)
$shape.Height = 1045

# ---------------------------------------------------------------------------
# 5. Misc. view state: the active cell ends up on G11 and the workbook
#    window is repositioned slightly on screen.
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 2700
$win.Top = 140
$ws.Range("G11").Select() | Out-Null

